$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "subcategories" (index 2 / sheet2.xml)
#   - insert 3 new rows after row 9 (Sports, Clothes, Household Items)
#   - insert 1 new row after row 17 (Tvs)
#   - fill column A ("Sub id") for every data row
#   - fill column C for the 4 newly inserted rows
# ---------------------------------------------------------------
$wsSub = $wb.Worksheets.Item(2)

$wsSub.Range("A10:A12").EntireRow.Insert()
$wsSub.Range("A18").EntireRow.Insert()

$wsSub.Range("C10").Value = "Sports"
$wsSub.Range("C11").Value = "Clothes"
$wsSub.Range("C12").Value = "Household Items"
$wsSub.Range("C18").Value = "Tvs "

$wsSub.Range("A2").Value = "Sid01"
$wsSub.Range("A3").Value = "Sid02"
$wsSub.Range("A4").Value = "Sid03"
$wsSub.Range("A5").Value = "Sid04"
$wsSub.Range("A6").Value = "Sid05"
$wsSub.Range("A7").Value = "Sid06"
$wsSub.Range("A8").Value = "Sid07"
$wsSub.Range("A9").Value = "Sid08"
$wsSub.Range("A10").Value = "Sid09"
$wsSub.Range("A11").Value = "Sid10"
$wsSub.Range("A12").Value = "Sid11"
$wsSub.Range("A13").Value = "Sid"
$wsSub.Range("A14").Value = "Sid13"
$wsSub.Range("A15").Value = "Sid"
$wsSub.Range("A16").Value = "Sid"
$wsSub.Range("A17").Value = "Sid"
$wsSub.Range("A18").Value = "Sid16"
$wsSub.Range("A19").Value = "Sid"
$wsSub.Range("A20").Value = "Sid"
$wsSub.Range("A21").Value = "Sid"
$wsSub.Range("A22").Value = "Sid"
$wsSub.Range("A23").Value = "Sid"
$wsSub.Range("A24").Value = "Sid"
$wsSub.Range("A25").Value = "Sid"
$wsSub.Range("A26").Value = "Sid"
$wsSub.Range("A27").Value = "Sid"
$wsSub.Range("A28").Value = "Sid"
$wsSub.Range("A29").Value = "Sid"
$wsSub.Range("A30").Value = "Sid"

$wsSub.Range("A11").Select()

# ---------------------------------------------------------------
# Sheet "vendor type" (index 3 / sheet3.xml)
#   - add a new row with a stray "_+" value
# ---------------------------------------------------------------
$wsVtype = $wb.Worksheets.Item(3)
$wsVtype.Range("B9").Value = "_+"
$wsVtype.Activate()
$excel.ActiveWindow.Zoom = 130
$wsVtype.Range("G10").Select()

# ---------------------------------------------------------------
# Sheet "products" (index 5 / sheet5.xml)
#   - fill in vendor (D) and subcategory (F) columns
#   - center align the C:F data block (creates the new shared style)
# ---------------------------------------------------------------
$wsProd = $wb.Worksheets.Item(5)

$wsProd.Range("D2").Value = "v7"
$wsProd.Range("F2").Value = "Sid01"

$wsProd.Range("D3").Value = "v2"
$wsProd.Range("F3").Value = "Sid16"

$wsProd.Range("D4").Value = "v6"
$wsProd.Range("F4").Value = "Sid11"

$wsProd.Range("D5").Value = "v5"
$wsProd.Range("F5").Value = "Sid06"

$wsProd.Range("D6").Value = "v8"
$wsProd.Range("F6").Value = "Sid13"

$wsProd.Range("D7").Value = "v8"
$wsProd.Range("F7").Value = "Sid13"

$wsProd.Range("D8").Value = "v7"
$wsProd.Range("F8").Value = "Sid02"

$wsProd.Range("D9").Value = "v7"
$wsProd.Range("F9").Value = "Sid08"

$wsProd.Range("D10").Value = "v7"
$wsProd.Range("F10").Value = "Sid09"

$wsProd.Range("D11").Value = "v8"
$wsProd.Range("F11").Value = "Sid10"

$wsProd.Range("C1:F11").HorizontalAlignment = -4108

$wsProd.Activate()
$wsProd.Range("D12").Select()
